$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (incl. new column H) ---
$ws.Columns.Item(1).ColumnWidth = 29.42
$ws.Columns.Item(2).ColumnWidth = 29.59
$ws.Columns.Item(3).ColumnWidth = 24.76
$ws.Columns.Item(4).ColumnWidth = 22.76
$ws.Columns.Item(5).ColumnWidth = 17.92
$ws.Columns.Item(6).ColumnWidth = 19.76
$ws.Columns.Item(7).ColumnWidth = 20.92
$ws.Columns.Item(8).ColumnWidth = 20.25

# --- Content: add the new "isVerified" column ---
$ws.Range("H1").Value = "{d.i18n.isVerified}"
$ws.Range("H2").Value = "{d.contacts[i].isVerified}"
$ws.Range("H3").Value = "{d.contacts[i+1].isVerified}"

# --- Content: add 7 new blank rows (4-10) under the table ---
for ($r = 4; $r -le 10; $r++) {
  $ws.Range("A" + $r + ":H" + $r).Value = ""
}

# Start every cell in the table from a clean, identical base so that
# cells which end up with the same formatting collapse to one shared style.
$full = $ws.Range("A1:H10")
$full.ClearFormats()

# ---------- Row 1: header ----------
$row1 = $ws.Range("A1:H1")
$row1.Font.Name = "Arial"
$row1.Font.Size = 10
$row1.Font.Bold = $true
$row1.Font.Color = 0
$row1.Interior.Pattern = 1
$row1.Interior.Color = 11711407
$row1.NumberFormat = "General"
$row1.VerticalAlignment = -4107
foreach ($edge in 7,8,9,10,11,12) {
  $row1.Borders.Item($edge).LineStyle = 1
  $row1.Borders.Item($edge).Weight = 2
  $row1.Borders.Item($edge).Color = 9737364
}
# bottom edge (separator above the data rows) is a darker gray
$row1.Borders.Item(9).Color = 3158064
$row1.Borders.Item(12).Color = 3158064

# ---------- Row 2 ----------
$row2 = $ws.Range("A2:H2")
$row2.Font.Name = "Arial"
$row2.Font.Size = 10
$row2.Font.Bold = $false
$row2.Font.Color = 0
$row2.Interior.Pattern = 0
$row2.NumberFormat = "General"
$row2.VerticalAlignment = -4107
foreach ($edge in 7,8,9,10,11) {
  $row2.Borders.Item($edge).LineStyle = 1
  $row2.Borders.Item($edge).Weight = 2
  $row2.Borders.Item($edge).Color = 9737364
}
$row2.Borders.Item(8).Color = 3158064

# ---------- Row 3 ----------
$row3 = $ws.Range("A3:H3")
$row3.Font.Name = "Arial"
$row3.Font.Size = 10
$row3.Font.Bold = $false
$row3.Font.Color = 0
$row3.Interior.Pattern = 0
$row3.NumberFormat = "General"
$row3.VerticalAlignment = -4107
foreach ($edge in 7,8,9,10,11,12) {
  $row3.Borders.Item($edge).LineStyle = 1
  $row3.Borders.Item($edge).Weight = 2
  $row3.Borders.Item($edge).Color = 9737364
}

# ---------- Row 4 ----------
$row4 = $ws.Range("A4:H4")
$row4.Font.Name = "Arial"
$row4.Font.Size = 10
$row4.Font.Bold = $false
$row4.Font.Color = 0
$row4.Interior.Pattern = 0
$row4.NumberFormat = "General"
$row4.VerticalAlignment = -4107
foreach ($edge in 7,8,9,10,11,12) {
  $row4.Borders.Item($edge).LineStyle = 1
  $row4.Borders.Item($edge).Weight = 2
  $row4.Borders.Item($edge).Color = 10132122
}
$row4.Borders.Item(8).Color = 9737364

# ---------- Rows 5-10 ----------
$rows510 = $ws.Range("A5:H10")
$rows510.Font.Name = "Arial"
$rows510.Font.Size = 10
$rows510.Font.Bold = $false
$rows510.Font.Color = 0
$rows510.Interior.Pattern = 0
$rows510.NumberFormat = "General"
$rows510.VerticalAlignment = -4107
foreach ($edge in 7,8,9,10,11,12) {
  $rows510.Borders.Item($edge).LineStyle = 1
  $rows510.Borders.Item($edge).Weight = 2
  $rows510.Borders.Item($edge).Color = 10132122
}

Write-Host "rows 1-10 formatted"
